# Auto-generated edit script applying scheduled-runner value updates
# to the Lamia Leve Profits workbook, per sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("H40").Value = 7301.4
$ws.Range("I40").Value = 3327.9412
$ws.Range("J40").Value = 9079
$ws.Range("K40").Value = 3327.9412
$ws.Range("L40").Value = 9079
$ws.Range("M40").Value = -3152.9412
$ws.Range("N40").Value = -9429
$ws.Range("H69").Value = 15000
$ws.Range("I69").Value = 12000
$ws.Range("J69").Value = 18000
$ws.Range("K69").Value = 36000
$ws.Range("L69").Value = 54000
$ws.Range("M69").Value = -35126
$ws.Range("N69").Value = -55748
$ws.Range("H72").Value = 15000
$ws.Range("I72").Value = 12000
$ws.Range("J72").Value = 18000
$ws.Range("K72").Value = 108000
$ws.Range("L72").Value = 162000
$ws.Range("M72").Value = -103632
$ws.Range("N72").Value = -170736
$ws.Range("H86").Value = 6237.6665
$ws.Range("J86").Value = 5607.4287
$ws.Range("L86").Value = 5607.4287
$ws.Range("N86").Value = -7853.4287
$ws.Range("H89").Value = 6237.6665
$ws.Range("J89").Value = 5607.4287
$ws.Range("L89").Value = 28037.1435
$ws.Range("N89").Value = -39269.14350000001
$ws.Range("H101").Value = 408
$ws.Range("I101").Value = 168
$ws.Range("J101").Value = 888
$ws.Range("K101").Value = 504
$ws.Range("L101").Value = 2664
$ws.Range("M101").Value = 1118
$ws.Range("N101").Value = -5908
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("H106").Value = 10269.5
$ws.Range("I106").Value = 1492.7142
$ws.Range("K106").Value = 1492.7142
$ws.Range("M106").Value = -861.7141999999999
$ws.Range("H112").Value = 1894.0769
$ws.Range("I112").Value = 237
$ws.Range("J112").Value = 2195.3635
$ws.Range("K112").Value = 711
$ws.Range("L112").Value = 6586.0905
$ws.Range("M112").Value = 397
$ws.Range("N112").Value = -8802.0905
$ws.Range("H116").Value = 12598.9
$ws.Range("I116").Value = 10665.667
$ws.Range("K116").Value = 10665.667
$ws.Range("M116").Value = -7223.666999999999
$ws.Range("H138").Value = 2845.1333
$ws.Range("I138").Value = 1652.4117
$ws.Range("J138").Value = 3569.2856
$ws.Range("K138").Value = 4957.2351
$ws.Range("L138").Value = 10707.8568
$ws.Range("M138").Value = 182.7649000000001
$ws.Range("N138").Value = -20987.8568
$ws.Range("N3").ClearContents()
$ws.Range("N102").ClearContents()
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4144.4614
$ws.Range("J122").Value = 7249.75
$ws.Range("L122").Value = 21749.25
$ws.Range("N122").Value = -26649.25
$ws.Range("H132").Value = 4015.3125
$ws.Range("I132").Value = 1969.3334
$ws.Range("K132").Value = 5908.0002
$ws.Range("M132").Value = -3378.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3269.0908
$ws.Range("I20").Value = 2506.8
$ws.Range("K20").Value = 2506.8
$ws.Range("M20").Value = -2259.8
$ws.Range("H69").Value = 40000
$ws.Range("J69").Value = 40000
$ws.Range("L69").Value = 40000
$ws.Range("N69").Value = -41622
$ws.Range("H72").Value = 40000
$ws.Range("J72").Value = 40000
$ws.Range("L72").Value = 120000
$ws.Range("N72").Value = -128112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5699.3335
$ws.Range("I22").Value = 5149
$ws.Range("K22").Value = 5149
$ws.Range("M22").Value = -4799
$ws.Range("H58").Value = 9443.182000000001
$ws.Range("I58").Value = 2371
$ws.Range("K58").Value = 2371
$ws.Range("M58").Value = -2168
$ws.Range("H75").Value = 45000
$ws.Range("J75").Value = 45000
$ws.Range("L75").Value = 45000
$ws.Range("N75").Value = -46996
$ws.Range("H78").Value = 45000
$ws.Range("J78").Value = 45000
$ws.Range("L78").Value = 135000
$ws.Range("N78").Value = -144984
$ws.Range("H134").Value = 3416.45
$ws.Range("I134").Value = 2244.125
$ws.Range("K134").Value = 6732.375
$ws.Range("M134").Value = -4197.375
$ws.Range("H135").Value = 64998.5
$ws.Range("J135").Value = 64998.5
$ws.Range("L135").Value = 64998.5
$ws.Range("N135").Value = -75138.5
$ws.Range("H136").Value = 9443.182000000001
$ws.Range("I136").Value = 2371
$ws.Range("K136").Value = 7113
$ws.Range("M136").Value = -4563

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 5667.6665
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 601880
$ws.Range("I3").Value = 668311.1
$ws.Range("K3").Value = 668311.1
$ws.Range("M3").Value = -668195.1
$ws.Range("H24").Value = 25094
$ws.Range("J24").Value = 26801.777
$ws.Range("L24").Value = 26801.777
$ws.Range("N24").Value = -27147.777
$ws.Range("H70").Value = 16974.941
$ws.Range("I70").Value = 5325.1816
$ws.Range("K70").Value = 5325.1816
$ws.Range("M70").Value = -5055.1816
$ws.Range("H73").Value = 16974.941
$ws.Range("I73").Value = 5325.1816
$ws.Range("K73").Value = 5325.1816
$ws.Range("M73").Value = -4389.1816
$ws.Range("H113").Value = 3538.5
$ws.Range("I113").Value = 1861.6
$ws.Range("K113").Value = 1861.6
$ws.Range("M113").Value = 308.4000000000001
$ws.Range("H132").Value = 59838.75
$ws.Range("I132").Value = 108178.3
$ws.Range("K132").Value = 324534.9
$ws.Range("M132").Value = -322004.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4483.7144
$ws.Range("I16").Value = 2877.6
$ws.Range("K16").Value = 2877.6
$ws.Range("M16").Value = -2707.6
$ws.Range("H61").Value = 9471.299999999999
$ws.Range("I61").Value = 8634.223
$ws.Range("K61").Value = 8634.223
$ws.Range("M61").Value = -8432.223
$ws.Range("H106").Value = 35125.715
$ws.Range("J106").Value = 35125.715
$ws.Range("L106").Value = 35125.715
$ws.Range("N106").Value = -37649.715
$ws.Range("H113").Value = 9471.299999999999
$ws.Range("I113").Value = 8634.223
$ws.Range("K113").Value = 8634.223
$ws.Range("M113").Value = -6464.223
$ws.Range("H136").Value = 3591.724
$ws.Range("I136").Value = 1797.9166
$ws.Range("K136").Value = 5393.7498
$ws.Range("M136").Value = -2843.7498

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8013.2
$ws.Range("I81").Value = 3688.3333
$ws.Range("K81").Value = 7376.6666
$ws.Range("M81").Value = -6315.6666
$ws.Range("H84").Value = 8013.2
$ws.Range("I84").Value = 3688.3333
$ws.Range("K84").Value = 36883.333
$ws.Range("M84").Value = -31579.333
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("H122").Value = 12208
$ws.Range("I122").Value = 3867.75
$ws.Range("K122").Value = 11603.25
$ws.Range("M122").Value = -9153.25
$ws.Range("N105").ClearContents()
